# aggiunto tecnologia al file excel
# Add two new technologies (Photovoltaic, Wind_Turbine) to the
# "Existing_capacities" sheet, with explanatory comments, and leave that
# sheet as the active/selected one (as last touched by the author).

$wb = $excel.ActiveWorkbook
$wsCap = $wb.Worksheets.Item("Existing_capacities")

# --- New technology rows --------------------------------------------------
$wsCap.Range("A5").Value = "Photovoltaic"
$wsCap.Range("A6").Value = "Wind_Turbine"

# --- Comments explaining existing / new technology rows -------------------
$c1 = $wsCap.Range("A4")
$c1.AddComment("Bertoni, L. (Luca):`nTutti I tipi") | Out-Null

$c2 = $wsCap.Range("A6")
$c2.AddComment("Bertoni, L. (Luca):`nSia onshore che offshore") | Out-Null

# --- Leave "Existing_capacities" as the active sheet / selection ----------
$wsCap.Activate()
$wsCap.Range("E20").Select() | Out-Null
